# Insert a new weekly price record for Cilantro (Terminal La Palmera de La
# Serena) as row 28, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28 - this shifts existing rows 28..125 down
# to 29..126 and extends the used range to A1:R126.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly observation.
$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44623
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 100112040
$ws.Cells.Item(28, 7).Value = "Cilantro"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 2200
$ws.Cells.Item(28, 11).Value = 2500
$ws.Cells.Item(28, 12).Value = 3000
$ws.Cells.Item(28, 13).Value = 2750
$ws.Cells.Item(28, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(28, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(28, 16).Value = 1833
$ws.Cells.Item(28, 17).Value = 1.5
$ws.Cells.Item(28, 18).Value = "Hortaliza"
